$p = $ppt.ActivePresentation
$s15 = $p.Slides.Add(15, 12)
$s16 = $p.Slides.Add(16, 12)

$s16.Background.Fill.Solid()
$s16.Background.Fill.ForeColor.ObjectThemeColor = 13
Write-Output "done"
